$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 'Natalie O''Rourke'
$ws.Cells.Item(3, 3).Value = 'Justyna Swierz'
$ws.Cells.Item(9, 3).Value = 'Owen Nevaril'
$ws.Cells.Item(10, 3).Value = 'Alexander Betz'
$ws.Cells.Item(12, 3).Value = 'Tran Nguyen'
$ws.Cells.Item(13, 3).Value = 'Natalie O''Rourke'
$ws.Cells.Item(17, 3).Value = 'Natalie O''Rourke'
$ws.Cells.Item(18, 3).Value = 'Justyna Swierz'
$ws.Cells.Item(19, 3).Value = 'Tran Nguyen'
$ws.Cells.Item(22, 3).Value = 'Justyna Swierz'
$ws.Cells.Item(24, 3).Value = 'Tran Nguyen'
$ws.Cells.Item(25, 3).Value = 'Vivian'
$ws.Cells.Item(26, 3).Value = 'Noah Williams'
$ws.Cells.Item(27, 3).Value = 'Owen Nevaril'
$ws.Cells.Item(28, 3).Value = 'Aspen'
$ws.Cells.Item(29, 3).Value = 'Tran Nguyen'
$ws.Cells.Item(30, 3).Value = 'Justyna Swierz'
$ws.Cells.Item(31, 3).Value = 'Aspen'
$ws.Cells.Item(33, 3).Value = 'Justyna Swierz'
$ws.Cells.Item(37, 3).Value = 'Alexander Betz'
$ws.Cells.Item(38, 3).Value = 'Tran Nguyen'
$ws.Cells.Item(39, 3).Value = 'Justyna Swierz'
$ws.Cells.Item(42, 3).Value = 'Tran Nguyen'
$ws.Cells.Item(44, 3).Value = 'Vivian'
$ws.Cells.Item(45, 3).Value = 'Justyna Swierz'
$ws.Cells.Item(46, 3).Value = 'Tran Nguyen'
$ws.Cells.Item(47, 3).Value = 'Vivian'
$ws.Cells.Item(49, 3).Value = 'Noah Williams'
$ws.Cells.Item(50, 3).Value = 'Natalie O''Rourke'
$ws.Cells.Item(54, 3).Value = 'Tran Nguyen'
$ws.Cells.Item(56, 3).Value = 'Vivian'
$ws.Cells.Item(57, 3).Value = 'Owen Nevaril'
$ws.Cells.Item(58, 3).Value = 'Owen Nevaril'
$ws.Cells.Item(59, 3).Value = 'Noah Williams'
$ws.Cells.Item(64, 3).Value = 'Alexander Betz'
$ws.Cells.Item(65, 3).Value = 'Noah Williams'
$ws.Cells.Item(66, 3).Value = 'Natalie O''Rourke'
$ws.Cells.Item(68, 3).Value = 'Alexander Betz'
$ws.Cells.Item(69, 3).Value = 'Aspen'
$ws.Cells.Item(72, 3).Value = 'Owen Nevaril'
$ws.Cells.Item(74, 3).Value = 'Noah Williams'
$ws.Cells.Item(75, 3).Value = 'Tran Nguyen'
$ws.Cells.Item(76, 3).Value = 'Natalie O''Rourke'
$ws.Cells.Item(77, 3).Value = 'Vivian'
$ws.Cells.Item(78, 3).Value = 'Vivian'
$ws.Cells.Item(79, 3).Value = 'Alexander Betz'
$ws.Cells.Item(81, 3).Value = 'Noah Williams'
$ws.Cells.Item(82, 3).Value = 'Natalie O''Rourke'
$ws.Cells.Item(83, 3).Value = 'Justyna Swierz'
$ws.Cells.Item(86, 3).Value = 'Owen Nevaril'
$ws.Cells.Item(89, 3).Value = 'Natalie O''Rourke'
$ws.Cells.Item(90, 3).Value = 'Tran Nguyen'
$ws.Cells.Item(91, 3).Value = 'Noah Williams'
$ws.Cells.Item(92, 3).Value = 'Alexander Betz'
$ws.Cells.Item(93, 3).Value = 'Natalie O''Rourke'
$ws.Cells.Item(94, 3).Value = 'Natalie O''Rourke'
$ws.Cells.Item(95, 3).Value = 'Aspen'
$ws.Cells.Item(97, 3).Value = 'Aspen'
$ws.Cells.Item(98, 3).Value = 'Justyna Swierz'
$ws.Cells.Item(99, 3).Value = 'Natalie O''Rourke'
$ws.Cells.Item(100, 3).Value = 'Alexander Betz'
$ws.Cells.Item(101, 3).Value = 'Aspen'
$ws.Cells.Item(102, 3).Value = 'Owen Nevaril'
$ws.Cells.Item(103, 3).Value = 'Tran Nguyen'
$ws.Cells.Item(104, 3).Value = 'Vivian'
$ws.Cells.Item(105, 3).Value = 'Natalie O''Rourke'
$ws.Cells.Item(107, 3).Value = 'Tran Nguyen'
$ws.Cells.Item(112, 3).Value = 'Aspen'
$ws.Cells.Item(114, 3).Value = 'Aspen'
$ws.Cells.Item(116, 3).Value = 'Aspen'
$ws.Cells.Item(117, 3).Value = 'Natalie O''Rourke'
$ws.Cells.Item(118, 3).Value = 'Justyna Swierz'
$ws.Cells.Item(119, 3).Value = 'Vivian'
$ws.Cells.Item(120, 3).Value = 'Tran Nguyen'
$ws.Cells.Item(121, 3).Value = 'Vivian'
$ws.Cells.Item(122, 3).Value = 'Tran Nguyen'
$ws.Cells.Item(125, 3).Value = 'Noah Williams'
$ws.Cells.Item(128, 3).Value = 'Vivian'
$ws.Cells.Item(129, 3).Value = 'Natalie O''Rourke'
$ws.Cells.Item(130, 3).Value = 'Aspen'
$ws.Cells.Item(131, 3).Value = 'Noah Williams'
$ws.Cells.Item(134, 3).Value = 'Justyna Swierz'
$ws.Cells.Item(135, 3).Value = 'Owen Nevaril'
$ws.Cells.Item(138, 3).Value = 'Alexander Betz'
$ws.Cells.Item(139, 3).Value = 'Noah Williams'
$ws.Cells.Item(140, 3).Value = 'Aspen'
$ws.Cells.Item(141, 3).Value = 'Owen Nevaril'
$ws.Cells.Item(142, 3).Value = 'Owen Nevaril'
$ws.Cells.Item(143, 3).Value = 'Alexander Betz'
$ws.Cells.Item(145, 3).Value = 'Natalie O''Rourke'
$ws.Cells.Item(146, 3).Value = 'Justyna Swierz'
$ws.Cells.Item(147, 3).Value = 'Owen Nevaril'
$ws.Cells.Item(148, 3).Value = 'Aspen'
$ws.Cells.Item(149, 3).Value = 'Alexander Betz'
$ws.Cells.Item(150, 3).Value = 'Owen Nevaril'
$ws.Cells.Item(151, 3).Value = 'Tran Nguyen'
$ws.Cells.Item(152, 3).Value = 'Owen Nevaril'
$ws.Cells.Item(155, 3).Value = 'Alexander Betz'
